# Adds a new column E ("Fonte_CX 6200 Switch Series") to the 'Resultados' sheet,
# filling in page-reference source values for each requirement row, and
# updates several requirement answers (column C) with corrected/expanded text.
# Row 69 additionally flips from an unanswered (yellow) requirement to an
# answered (green) one with new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New column E header (copy formatting from D1, the last header cell) ----
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E1").Value = "Fonte_CX 6200 Switch Series"

# ---- Column E (Fonte) values for each data row ----
$eVals = @{}
$eVals[2] = "None"
$eVals[3] = "Page 25"
$eVals[4] = "Page 25"
$eVals[5] = "Page 25"
$eVals[6] = "Page 25"
$eVals[7] = "Page 25"
$eVals[8] = "Page 25"
$eVals[9] = "Page 25"
$eVals[10] = "Page 25"
$eVals[11] = "Page 25"
$eVals[12] = "Page 25"
$eVals[13] = "Page 25"
$eVals[14] = "Page 3"
$eVals[15] = "Page 4"
$eVals[16] = "Page 28"
$eVals[17] = "Page 25"
$eVals[18] = "Page 25"
$eVals[19] = "None"
$eVals[20] = "Page 6"
$eVals[21] = "Page 6"
$eVals[22] = "Page 6"
$eVals[23] = "Page 6"
$eVals[24] = "Page 6"
$eVals[25] = "Page 4"
$eVals[26] = "Page 6"
$eVals[27] = "Page 6"
$eVals[28] = "Page 5"
$eVals[29] = "Page 22"
$eVals[30] = "Page 6"
$eVals[31] = "Page 6"
$eVals[32] = "Page 6"
$eVals[33] = "Page 6"
$eVals[34] = "Page 6"
$eVals[35] = "Page 22"
$eVals[36] = "Page 6"
$eVals[37] = "Page 6"
$eVals[38] = "Page 6"
$eVals[39] = "Page 6"
$eVals[40] = "Page 22"
$eVals[41] = "Page 4"
$eVals[42] = "Page 7"
$eVals[43] = "Page 5"
$eVals[44] = "Page 5"
$eVals[45] = "Page 5"
$eVals[46] = "Page 5"
$eVals[47] = "None"
$eVals[48] = "Page 7"
$eVals[49] = "Page 7"
$eVals[50] = "Page 3"
$eVals[51] = "Page 7"
$eVals[52] = "Page 7"
$eVals[53] = "Page 7"
$eVals[54] = "Page 7"
$eVals[55] = "Page 7"
$eVals[56] = "Page 7"
$eVals[57] = "Page 5"
$eVals[58] = "Page 5"
$eVals[59] = "Page 5"
$eVals[60] = "None"
$eVals[61] = "Page 5"
$eVals[62] = "Page 5"
$eVals[63] = "Page 7"
$eVals[64] = "Page 5"
$eVals[65] = "Page 5"
$eVals[66] = "Page 5"
$eVals[67] = "Page 8"
$eVals[68] = "Page 3"
$eVals[69] = "Page 3"

# ---- Updated column C (Resposta) values ----
$cVals = @{}
$cVals[3] = "12 ports SmartRate 100M/1G/2.5G/5G BaseT"
$cVals[4] = "4x 100M/1G/10G SFP ports (2x LRM"
$cVals[5] = "4x 100M/1G/10G SFP ports (2x LRM/MACSec 256)"
$cVals[14] = "Support for up to 8 switches (or members) in a stack via chain or ring topology"
$cVals[15] = "1U"
$cVals[16] = "Mounts in an EIA-standard 19 in. Telco rack or equipment cabinet."
$cVals[18] = "2 field-replaceable, hot-swappable power supply slots `n1 minimum power supply required (ordered separately)"
$cVals[21] = "IEEE 802.1v protocol VLANs"
$cVals[22] = "Bridge Protocol Data Unit (BPDU) tunneling"
$cVals[23] = "Jumbo packet support improves the performance `nof large data transfers; supports frame size of up to `n9198 bytes"
$cVals[24] = "Port mirroring duplicates port traffic (ingress `nand egress) to a monitoring port; supports `n4 mirroring groups"
$cVals[27] = "MVRP allows automatic learning and dynamic `nassignment of VLANs"
$cVals[29] = "LLDP-MED (Media Endpoint Discovery) defines a `nstandard extension of LLDP"
$cVals[32] = "MVRP allows automatic learning and dynamic `nassignment of VLANs"
$cVals[34] = "VXLAN encapsulation tunneling protocol for `noverlay network that enables a more scalable virtual `nnetwork deployment"
$cVals[36] = "Static IP routing"
$cVals[38] = "OSPFv3 for IPv6 routing"
$cVals[39] = "DHCP server"
$cVals[40] = "2,048 `n1,024"
$cVals[41] = "Packet storm protection"
$cVals[42] = "ICMP throttling"
$cVals[43] = "Strict priority (SP) queuing and Deficit Weighted`nRound Robin (DWRR)"
$cVals[44] = "Traffic prioritization (IEEE 802.1p) for`nreal-time classification"
$cVals[46] = "Class of Service (CoS) sets the IEEE 802.1p priority `ntag based on IP address, IP Type of Service (ToS), `nLayer 3 protocol, TCP/UDP port number, source port,`nand DiffServ"
$cVals[48] = "Access control list (ACL) support for both IPv4 and`nIPv6"
$cVals[49] = "ACLs also provide filtering based on the IP field, `nsource/destination IP address/subnet, and `nsource/destination TCP/UDP port number on a `nper-VLAN or per-port basis"
$cVals[50] = "Role-based microsegmentation"
$cVals[51] = "Uses `nan IEEE 802.1X supplicant on the client"
$cVals[54] = "Concurrent IEEE 802.1X, Web, and MAC authentication `nschemes per switch port accepts up to 32 sessions of `nIEEE 802.1X, Web, and MAC authentications"
$cVals[55] = "Terminal Access Controller Access-Control System`n(TACACS+)"
$cVals[56] = "Switch CPU protection provides automatic protection`nagainst malicious network traffic trying to shut down`nthe switch"
$cVals[58] = "Dual flash images provides independent primary`nand secondary operating system files for backup`nwhile upgrading"
$cVals[59] = "Multiple configuration files can be stored to a`nflash image"
$cVals[61] = "Ingress and egress port monitoring enable more`nefficient network problem solving"
$cVals[62] = "sFlow® (RFC 3176)"
$cVals[63] = "allowing secure access to the browser-based `nmanagement GUI in the switch"
$cVals[64] = "local and remote`nsyslog capabilities"
$cVals[65] = "Secure File Transfer Protocol (SFTP)"
$cVals[66] = "Supports SNMP (v2c/v3)"
$cVals[68] = "An easy-to-use mobile app simplifies connecting, `nstacking and managing HPE Aruba Networking CX `n6200 switches for any size project."
$cVals[69] = "Dynamic Segmentation provides scale and flexibility in `nnetwork design by allowing the stretching of VLANs and `nsubnets across the entire network with a VXLAN-based `ndistributed overlay fabric."

foreach ($row in $eVals.Keys) {
    $ws.Cells.Item($row, 5).Value = $eVals[$row]
}

foreach ($row in $cVals.Keys) {
    $ws.Cells.Item($row, 3).Value = $cVals[$row]
}

# ---- Row 69 special case: requirement now met (was 'None'/yellow) ----
$ws.Range("C69").Interior.Color = 32768
$ws.Range("D69").Value = "green"

